$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.284.66"
$ws.Range("E2").Value = "  +2.82%  "
$ws.Range("D3").Value = "1.895.43"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -1.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.23"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5145"
$ws.Range("E7").Value = "  +0.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3916"
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08425"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.46"
$ws.Range("E10").Value = "  +1.62%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.243"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.895.08"
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.66"
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.06"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06744"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.81"
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.014"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("D23").Value = "29.291.34"
$ws.Range("E23").Value = "  +2.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.13"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.213"
$ws.Range("E25").Value = "  -2.50%  "
$ws.Range("D26").Value = "2.114.45"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.16"
$ws.Range("E27").Value = "  -1.22%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.90"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.436"
$ws.Range("E29").Value = "  +2.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.93"
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("E31").Value = "  +0.76%  "
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.127"
$ws.Range("E33").Value = "  +5.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.648"
$ws.Range("E34").Value = "  +0.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02477"
$ws.Range("E35").Value = "  +1.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06536"
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("E37").Value = "  +0.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2193"
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.226"
$ws.Range("E39").Value = "  +2.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.125"
$ws.Range("E40").Value = "  +1.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6505"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.232"
$ws.Range("E42").Value = "  -2.79%  "
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6058"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.09"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("E47").Value = "  +1.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.229"
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.16"
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("E50").Value = "  -2.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.63"
